$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text changes (E2:F3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("K2").Value = "2016-08-25 19:05:49"
$wsZh.Range("K3").Value = "2016-08-25 19:05:49"
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ad62a3a3f1cf178a58fa6f86b35b7016d7187ea/e2e/fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md", "", "", "fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md")
$wsZh.Range("I2").Value = "fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ad62a3a3f1cf178a58fa6f86b35b7016d7187ea/e2e/fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md", "", "", "fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ad62a3a3f1cf178a58fa6f86b35b7016d7187ea/e2e/ffffc7cc1db2-99b0-4511-bbb1-1592252a3feb.md", "", "", "ffffc7cc1db2-99b0-4511-bbb1-1592252a3feb.md")
$wsZh.Range("I3").Value = "fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ad62a3a3f1cf178a58fa6f86b35b7016d7187ea/e2e/fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md", "", "", "fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md")
$wsZh.Range("J2").Value = "fa547c0d-cd1c-4e34-8711-6f34bda7bb82.dec929f9212e255ca06cd666eafdecb727f7f07f.zh-cn.xlf"
$wsZh.Range("J3").Value = "fa547c0d-cd1c-4e34-8711-6f34bda7bb82.dec929f9212e255ca06cd666eafdecb727f7f07f.zh-cn.xlf"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("K2").Value = "2016-08-25 19:05:55"
$wsDe.Range("K3").Value = "2016-08-25 19:05:55"
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ad62a3a3f1cf178a58fa6f86b35b7016d7187ea/e2e/fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md", "", "", "fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md")
$wsDe.Range("I2").Value = "fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ad62a3a3f1cf178a58fa6f86b35b7016d7187ea/e2e/fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md", "", "", "fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ad62a3a3f1cf178a58fa6f86b35b7016d7187ea/e2e/ffffc7cc1db2-99b0-4511-bbb1-1592252a3feb.md", "", "", "ffffc7cc1db2-99b0-4511-bbb1-1592252a3feb.md")
$wsDe.Range("I3").Value = "fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ad62a3a3f1cf178a58fa6f86b35b7016d7187ea/e2e/fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md", "", "", "fa547c0d-cd1c-4e34-8711-6f34bda7bb82.md")
$wsDe.Range("J2").Value = "fa547c0d-cd1c-4e34-8711-6f34bda7bb82.dec929f9212e255ca06cd666eafdecb727f7f07f.de-de.xlf"
$wsDe.Range("J3").Value = "fa547c0d-cd1c-4e34-8711-6f34bda7bb82.dec929f9212e255ca06cd666eafdecb727f7f07f.de-de.xlf"

# --- Column widths to mirror Excel's width recalculation after content changes ---
$wsOverview.Range("E:E").ColumnWidth = 29.9777047293527
$wsOverview.Range("F:F").ColumnWidth = 29.9777047293527
$wsZh.Range("C:C").ColumnWidth = 29.9777047293527
$wsZh.Range("I:I").ColumnWidth = 40
$wsZh.Range("J:J").ColumnWidth = 40
$wsDe.Range("C:C").ColumnWidth = 29.9777047293527
$wsDe.Range("I:I").ColumnWidth = 40
$wsDe.Range("J:J").ColumnWidth = 40
